$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark from the end of the document (just before
#    "How should the predictor maintain...") to the very start of the
#    document (just before "Summary: This is a rough draft...").
# ---------------------------------------------------------------------------

# Remove the existing bookmark near the end of the document.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Re-create it at document position 0. Adding a zero-length bookmark
# directly at position 0 is unreliable, so insert a throwaway character,
# anchor the (non-zero-offset) bookmark right after it, then delete the
# throwaway character again -- the bookmark collapses back to a true
# zero-width bookmark at the very start of the document without splitting
# the first run.
$startRng = $d.Range(0, 0)
$startRng.InsertBefore("X")
$afterPlaceholder = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $afterPlaceholder)
$placeholderRng = $d.Range(0, 1)
$placeholderRng.Text = ""

# ---------------------------------------------------------------------------
# 2. Merge the "This function will " / "use the provided information..."
#    runs (comic scheduler's Description) into a single run.
# ---------------------------------------------------------------------------

$descRng = $d.Content
$null = $descRng.Find.Execute("This function will use the provided information to update the list of comics to be checked when.  The schedule of new comics to be checked should be maintained either in a file or in the database if we decide to have one.")
$descBoundary = $descRng.Start + ("This function will ").Length
$descCharRng = $d.Range($descBoundary - 1, $descBoundary)
$descChar = $descCharRng.Text
$descCharRng.Text = ""
$descInsRng = $d.Range($descBoundary - 1, $descBoundary - 1)
$descInsRng.InsertAfter($descChar)

# ---------------------------------------------------------------------------
# 3. Merge the "Inputs:" / " " runs (scheduler's Inputs row) into a single
#    run "Inputs: ".
# ---------------------------------------------------------------------------

$inputsRng = $d.Content
$null = $inputsRng.Find.Execute("Inputs: none")
$inputsBoundary = $inputsRng.Start + ("Inputs:").Length
$inputsCharRng = $d.Range($inputsBoundary - 1, $inputsBoundary)
$inputsChar = $inputsCharRng.Text
$inputsCharRng.Text = ""
$inputsInsRng = $d.Range($inputsBoundary - 1, $inputsBoundary - 1)
$inputsInsRng.InsertAfter($inputsChar)

Write-Output "edits applied"
